# Adds two more years (2019, 2020) of data to the "Renewable energy share"
# table: one new column per year in the header row, the percentage row and
# the hydropower-production row, mirroring the existing 2007-2018 columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (year headers): P4 = 2019, Q4 = 2020 -----------------------
# Clone formatting from the last existing year column (O4) so the new
# cells pick up the same font/border/number format.
$ws.Cells.Item(4, 15).Copy()
$ws.Range($ws.Cells.Item(4, 16), $ws.Cells.Item(4, 17)).PasteSpecial(-4122)
$ws.Cells.Item(4, 16).Value = 2019
$ws.Cells.Item(4, 17).Value = 2020

# --- Row 5 (renewable energy share, %): P5 = 35.67, Q5 left blank -----
# Clone formatting from E5, which uses the percentage-style number format
# (style also used by the sibling cells on this row).
$ws.Cells.Item(5, 5).Copy()
$ws.Range($ws.Cells.Item(5, 16), $ws.Cells.Item(5, 17)).PasteSpecial(-4122)
$ws.Cells.Item(5, 16).Value = 35.67

# --- Row 6 (hydropower production, mln kWh): P6 = 13859.3, Q6 = 13979.1
# Clone formatting from the last existing year column (O6).
$ws.Cells.Item(6, 15).Copy()
$ws.Range($ws.Cells.Item(6, 16), $ws.Cells.Item(6, 17)).PasteSpecial(-4122)
$ws.Cells.Item(6, 16).Value = 13859.3
$ws.Cells.Item(6, 17).Value = 13979.1

# Clear the clipboard marquee left behind by Copy().
$excel.CutCopyMode = 0

# The new data makes the used range grow from A1:O7 to A1:Q7 and leaves
# the active selection on P9, matching the authored workbook.
$ws.Range("P9").Select()
